# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" / clrScheme "Office"      (currently unused, only
#                             wired to the Notes Master, which this PowerPoint host does not
#                             expose a live ThemeColorScheme for)
#   ppt/theme/theme2.xml  -> "Integral" / clrScheme "Red Violet"     (the theme actually in use
#                             by the one Slide Master / all slides & layouts)
#
# The authored edit swaps the two themes' content: the live theme (theme2.xml) becomes the
# plain "Office Theme" colours, and the dormant theme (theme1.xml) becomes "Integral". The
# font scheme and format scheme (fills/lines/effects) are identical between the two themes
# already, so only the 12 colour-scheme slots actually change.
#
# This host's object model only exposes a writable ThemeColorScheme through a Slide (it maps
# straight onto the single live theme part, theme2.xml) so we drive the swap through that:
# every slide shares the one theme, so touching slide 1's scheme updates the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# msoThemeColorDark1=1, Light1=2, Dark2=3, Light2=4, Accent1..6=5..10, Hyperlink=11, FollowedHyperlink=12
# Target ("Office Theme") colours, written as plain RRGGBB hex strings.
$officeThemeColors = @{
    1  = "000000"   # Dark 1
    2  = "FFFFFF"   # Light 1
    3  = "44546A"   # Dark 2
    4  = "E7E6E6"   # Light 2
    5  = "5B9BD5"   # Accent 1
    6  = "ED7D31"   # Accent 2
    7  = "A5A5A5"   # Accent 3
    8  = "FFC000"   # Accent 4
    9  = "4472C4"   # Accent 5
    10 = "70AD47"   # Accent 6
    11 = "0563C1"   # Hyperlink
    12 = "954F72"   # Followed Hyperlink
}

foreach ($idx in $officeThemeColors.Keys) {
    $hexStr = $officeThemeColors[$idx]
    $R = [Convert]::ToInt32($hexStr.Substring(0,2), 16)
    $G = [Convert]::ToInt32($hexStr.Substring(2,2), 16)
    $B = [Convert]::ToInt32($hexStr.Substring(4,2), 16)
    # VBA/COM RGB() longs pack as R + G*256 + B*65536
    $comRgb = $R + ($G * 256) + ($B * 65536)
    $tcs.Colors($idx).RGB = $comRgb
}
